# Rename both worksheets from "..._借记" to "..._原始凭证" and
# bump the service-number value (B12) on the first sheet from 2372 to 2882.

$wb = $excel.ActiveWorkbook

$oldName1 = "浙江杭州滨江中南乐游城店_借记"
$oldName2 = "浙江杭州三墩地铁站店_借记"
$newName1 = "浙江杭州滨江中南乐游城店_原始凭证"
$newName2 = "浙江杭州三墩地铁站店_原始凭证"

$ws1 = $wb.Worksheets($oldName1)
$ws2 = $wb.Worksheets($oldName2)

$ws1.Name = $newName1
$ws2.Name = $newName2

# Renaming the sheets updates every defined name that refers to a real
# range, but a defined name whose target was already an invalid
# reference (#REF!) loses its sheet qualifier instead of having it
# updated. Put the sheet-qualified #REF! back for "当前费率" on both
# scopes (workbook-level and the sheet-local copy on sheet 1).
foreach ($n in $wb.Names) {
    if ($n.Name -eq "当前费率") {
        $n.RefersTo = "=" + $newName2 + "!#REF!"
    }
    elseif ($n.Name -eq ($newName1 + "!当前费率") -or $n.Name -like "*!当前费率") {
        $n.RefersTo = "=" + $newName1 + "!#REF!"
    }
}

# Update the service number on the first sheet.
$ws1.Range("B12").Value = 2882
